$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 79; this shifts existing rows 79..276 down to 80..277
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with a new data record (same constant
# columns as the rest of the sheet, new Fecha + Volumen values)
$ws.Range("A79").Value = 3
$ws.Range("B79").Value = "Femacal de La Calera"
$ws.Range("C79").Value = "Coquimbo"
$ws.Range("D79").Value = 44622
$ws.Range("E79").Value = 5
$ws.Range("F79").Value = 100112039
$ws.Range("G79").Value = "Ciboulette"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 120
$ws.Range("K79").Value = 1500
$ws.Range("L79").Value = 1500
$ws.Range("M79").Value = 1500
$ws.Range("N79").Value = '$/docena de atados'
$ws.Range("O79").Value = "Provincia de Quillota"
$ws.Range("P79").Value = 500
$ws.Range("Q79").Value = 3
$ws.Range("R79").Value = "Hortaliza"
